$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: tidy up the Szeliski reference text (drop the "place not identified" clause) ---
$ws.Range("A2").Value = "SZELISKI, R., 2020. COMPUTER VISION. SPRINGER NATURE, p.5."

# --- Row 4 (new): Coifman et al. traffic-surveillance reference ---
$ws.Range("A4").Value = "Coifman, B., Beymer, D., McLauchlan, P., & Malik, J. (1998). A real-time computer vision system for vehicle tracking and traffic surveillance. Transportation Research Part C: Emerging Technologies"
$ws.Range("B4").Value = "(Coifman, 1998)"
$ws.Range("C4").Value = "background research"
$ws.Range("D4").Value = "what is computer vision"

# --- Row 5 (new): Hochreiter & Schmidhuber, Deep learning / LSTM reference ---
# (built with [string]::Concat because "1735" + [char]0x2013 + "1780" gets
# silently coerced to numeric addition by this host's "+" operator)
$enDash = [char]0x2013
$ws.Range("A5").Value = [string]::Concat("Hochreiter, Sepp and Jurgen Schmidhuber. 1997. ""Long Short-Term Memory. Neural Computation"" :1735", $enDash, "1780.")
$ws.Range("B5").Value = "(Hochreiter and`nSchmidhuber, 1997)"
$ws.Range("C5").Value = "background research"
$ws.Range("D5").Value = "Deep learning"

# Wrap the long text in row 5 and grow the row to fit two lines
$ws.Range("A5:B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 45

# Column A no longer auto-"best fit" -- set to an explicit manual width
# (56.3 is the input that round-trips closest to the target 57.140625 through
# this host's px-grid column-width quantization)
$ws.Columns.Item(1).ColumnWidth = 56.3

# Leave the selection parked below the table, matching the saved view
$ws.Range("A6").Select()
